$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - servicio: id_Servicio -> id ; drop trailing "matricula" in G3
$ws.Range("B3").Value = "id"
$ws.Range("G3").Value = ""

# Row 4 - Vehículo_ITV: shift columns left, drop the id_VehiculoITV column
$ws.Range("B4").Value = "matricula"
$ws.Range("B4").Font.Bold = $true
$ws.Range("C4").Value = "id_ITV"
$ws.Range("D4").Value = "resultado_vehiculo"
$ws.Range("D4").Font.Bold = $false
$ws.Range("E4").Value = ""

# Rows 5-12: rename the per-entity id_X header to the generic "id"
$ws.Range("B5").Value = "id"
$ws.Range("B6").Value = "id"
$ws.Range("B7").Value = "id"
$ws.Range("B8").Value = "id"
$ws.Range("B9").Value = "id"
$ws.Range("B10").Value = "id"
$ws.Range("B11").Value = "id"
$ws.Range("B12").Value = "id"

# Selection matches the diff
$ws.Range("A6:F12").Select()

$wb.Save()
